$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A10").Value = "GRT-USD"
